$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new shared string / task row ---
# Insert a new row at position 5 (existing rows 5.. shift down to 6..)
$ws.Rows("5:5").Insert()

# New row 5: "Research and get working example of multithreading" assigned to Mike
$ws.Range("A5").Value = "Research and get working example of multithreading"
$ws.Range("H6").Copy($ws.Range("B6"))
$ws.Range("H6").Copy($ws.Range("F6"))
$ws.Range("H6").Copy($ws.Range("B5"))
$ws.Range("H6").Copy($ws.Range("F5"))
$ws.Range("H6").Copy($ws.Range("H5"))
$ws.Range("H6").Clear()
$ws.Rows("5:5").RowHeight = 30

# Row 7 (previously row 6, "Creation of basic text UI"): add Andy markers in B and C
$ws.Range("H4").Copy($ws.Range("B7"))
$ws.Range("H4").Copy($ws.Range("C7"))

# Row 8 (previously row 7, "Ability to parse..."): add Yousef markers in B and C
$ws.Range("H3").Copy($ws.Range("B8"))
$ws.Range("H3").Copy($ws.Range("C8"))

# Row 9 (previously row 8, "Calculate a projected score..."): add Pat marker in B
$ws.Range("B4").Copy($ws.Range("B9"))

# Row 10 (previously row 9, "Creation of basic algorithm..."): add Pat markers in B and C
$ws.Range("B4").Copy($ws.Range("B10"))
$ws.Range("B4").Copy($ws.Range("C10"))

# Rows 3 and 4: duplicate the "Pat" marker from column B into column F
$ws.Range("B3").Copy($ws.Range("F3"))
$ws.Range("B4").Copy($ws.Range("F4"))

# Update the active selection to match the authored state
$ws.Range("B9").Select()
